# Atualizacao de bases das ligas: the betting-odds rows for several
# fixtures were re-ordered/re-matched (ids B + stats E:AD rotated among
# the rows that share the same match date), while the row index (A),
# Div (C) and Date (D) columns stay put. Apply the new values cell by
# cell for every affected row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Croatia 1NL")

# Row 9
$ws.Cells.Item(9, 2).Value = 6834732
$ws.Cells.Item(9, 5).Value = "Orijent"
$ws.Cells.Item(9, 6).Value = "NK Dubrava Zagreb"
$ws.Cells.Item(9, 7).Value = 1
$ws.Cells.Item(9, 8).Value = 1
$ws.Cells.Item(9, 9).Value = 1
$ws.Cells.Item(9, 10).Value = 1
$ws.Cells.Item(9, 11).Value = "D"
$ws.Cells.Item(9, 12).Value = 2.1
$ws.Cells.Item(9, 13).Value = 3.2
$ws.Cells.Item(9, 14).Value = 3.1
$ws.Cells.Item(9, 15).Value = 2.3
$ws.Cells.Item(9, 16).Value = 3.1
$ws.Cells.Item(9, 17).Value = 2.8
$ws.Cells.Item(9, 18).Value = -0.25
$ws.Cells.Item(9, 19).Value = 2.025
$ws.Cells.Item(9, 20).Value = 1.775
$ws.Cells.Item(9, 21).Value = 2.25
$ws.Cells.Item(9, 22).Value = 1.825
$ws.Cells.Item(9, 23).Value = 1.975
$ws.Cells.Item(9, 24).Value = -1
$ws.Cells.Item(9, 25).Value = 2.1
$ws.Cells.Item(9, 26).Value = -1
$ws.Cells.Item(9, 27).Value = -0.5
$ws.Cells.Item(9, 28).Value = 0.3875
$ws.Cells.Item(9, 29).Value = -0.5
$ws.Cells.Item(9, 30).Value = 0.4875

# Row 10
$ws.Cells.Item(10, 2).Value = 6834733
$ws.Cells.Item(10, 5).Value = "HNK Cibalia"
$ws.Cells.Item(10, 6).Value = "NK Croatia Zmijavci"
$ws.Cells.Item(10, 7).Value = 1
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 1
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = "H"
$ws.Cells.Item(10, 12).Value = 1.65
$ws.Cells.Item(10, 13).Value = 3.5
$ws.Cells.Item(10, 14).Value = 4.5
$ws.Cells.Item(10, 15).Value = 1.909
$ws.Cells.Item(10, 16).Value = 3.3
$ws.Cells.Item(10, 17).Value = 3.3
$ws.Cells.Item(10, 18).Value = -0.5
$ws.Cells.Item(10, 19).Value = 2
$ws.Cells.Item(10, 20).Value = 1.8
$ws.Cells.Item(10, 21).Value = 2.25
$ws.Cells.Item(10, 22).Value = 1.95
$ws.Cells.Item(10, 23).Value = 1.85
$ws.Cells.Item(10, 24).Value = 0.909
$ws.Cells.Item(10, 25).Value = -1
$ws.Cells.Item(10, 26).Value = -1
$ws.Cells.Item(10, 27).Value = 1
$ws.Cells.Item(10, 28).Value = -1
$ws.Cells.Item(10, 29).Value = -1
$ws.Cells.Item(10, 30).Value = 0.8500000000000001

# Row 11
$ws.Cells.Item(11, 2).Value = 6834729
$ws.Cells.Item(11, 5).Value = "NK Solin"
$ws.Cells.Item(11, 6).Value = "Bijelo Brdo"
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = "D"
$ws.Cells.Item(11, 12).Value = 2.1
$ws.Cells.Item(11, 13).Value = 3.2
$ws.Cells.Item(11, 14).Value = 3.1
$ws.Cells.Item(11, 15).Value = 2.05
$ws.Cells.Item(11, 16).Value = 3.25
$ws.Cells.Item(11, 17).Value = 3.25
$ws.Cells.Item(11, 18).Value = -0.25
$ws.Cells.Item(11, 19).Value = 1.8
$ws.Cells.Item(11, 20).Value = 2
$ws.Cells.Item(11, 21).Value = 2.5
$ws.Cells.Item(11, 22).Value = 2
$ws.Cells.Item(11, 23).Value = 1.8
$ws.Cells.Item(11, 24).Value = -1
$ws.Cells.Item(11, 25).Value = 2.25
$ws.Cells.Item(11, 26).Value = -1
$ws.Cells.Item(11, 27).Value = -0.5
$ws.Cells.Item(11, 28).Value = 0.5
$ws.Cells.Item(11, 29).Value = -1
$ws.Cells.Item(11, 30).Value = 0.8

# Row 126
$ws.Cells.Item(126, 2).Value = 6834841
$ws.Cells.Item(126, 5).Value = "Vukovar 91"
$ws.Cells.Item(126, 6).Value = "NK Jarun"
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = 1
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = "A"
$ws.Cells.Item(126, 12).Value = 1.5
$ws.Cells.Item(126, 13).Value = 4.2
$ws.Cells.Item(126, 14).Value = 4.75
$ws.Cells.Item(126, 15).Value = 1.4
$ws.Cells.Item(126, 16).Value = 4.75
$ws.Cells.Item(126, 17).Value = 5.5
$ws.Cells.Item(126, 18).Value = -1.25
$ws.Cells.Item(126, 19).Value = 1.95
$ws.Cells.Item(126, 20).Value = 1.85
$ws.Cells.Item(126, 21).Value = 2.75
$ws.Cells.Item(126, 22).Value = 1.95
$ws.Cells.Item(126, 23).Value = 1.85
$ws.Cells.Item(126, 24).Value = -1
$ws.Cells.Item(126, 25).Value = -1
$ws.Cells.Item(126, 26).Value = 4.5
$ws.Cells.Item(126, 27).Value = -1
$ws.Cells.Item(126, 28).Value = 0.8500000000000001
$ws.Cells.Item(126, 29).Value = -1
$ws.Cells.Item(126, 30).Value = 0.8500000000000001

# Row 127
$ws.Cells.Item(127, 2).Value = 6834840
$ws.Cells.Item(127, 5).Value = "NK Solin"
$ws.Cells.Item(127, 6).Value = "NK Dubrava Zagreb"
$ws.Cells.Item(127, 7).Value = 3
$ws.Cells.Item(127, 8).Value = 1
$ws.Cells.Item(127, 9).Value = 2
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 11).Value = "H"
$ws.Cells.Item(127, 12).Value = 2.5
$ws.Cells.Item(127, 13).Value = 3.2
$ws.Cells.Item(127, 14).Value = 2.5
$ws.Cells.Item(127, 15).Value = 2.5
$ws.Cells.Item(127, 16).Value = 3.2
$ws.Cells.Item(127, 17).Value = 2.5
$ws.Cells.Item(127, 18).Value = 0
$ws.Cells.Item(127, 19).Value = 1.9
$ws.Cells.Item(127, 20).Value = 1.9
$ws.Cells.Item(127, 21).Value = 2.25
$ws.Cells.Item(127, 22).Value = 1.825
$ws.Cells.Item(127, 23).Value = 1.975
$ws.Cells.Item(127, 24).Value = 1.5
$ws.Cells.Item(127, 25).Value = -1
$ws.Cells.Item(127, 26).Value = -1
$ws.Cells.Item(127, 27).Value = 0.8999999999999999
$ws.Cells.Item(127, 28).Value = -1
$ws.Cells.Item(127, 29).Value = 0.825
$ws.Cells.Item(127, 30).Value = -1

# Row 130
$ws.Cells.Item(130, 2).Value = 6834844
$ws.Cells.Item(130, 5).Value = "NK Jarun"
$ws.Cells.Item(130, 6).Value = "NK Solin"
$ws.Cells.Item(130, 7).Value = 2
$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 9).Value = 1
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 11).Value = "H"
$ws.Cells.Item(130, 12).Value = 1.727
$ws.Cells.Item(130, 13).Value = 3.6
$ws.Cells.Item(130, 14).Value = 4
$ws.Cells.Item(130, 15).Value = 2.05
$ws.Cells.Item(130, 16).Value = 3.4
$ws.Cells.Item(130, 17).Value = 3.1
$ws.Cells.Item(130, 18).Value = -0.25
$ws.Cells.Item(130, 19).Value = 1.85
$ws.Cells.Item(130, 20).Value = 1.95
$ws.Cells.Item(130, 21).Value = 2.5
$ws.Cells.Item(130, 22).Value = 1.925
$ws.Cells.Item(130, 23).Value = 1.875
$ws.Cells.Item(130, 24).Value = 1.05
$ws.Cells.Item(130, 25).Value = -1
$ws.Cells.Item(130, 26).Value = -1
$ws.Cells.Item(130, 27).Value = 0.8500000000000001
$ws.Cells.Item(130, 28).Value = -1
$ws.Cells.Item(130, 29).Value = -1
$ws.Cells.Item(130, 30).Value = 0.875

# Row 131
$ws.Cells.Item(131, 2).Value = 6834845
$ws.Cells.Item(131, 5).Value = "NK Dubrava Zagreb"
$ws.Cells.Item(131, 6).Value = "HNK Sibenik"
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 3
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 2
$ws.Cells.Item(131, 11).Value = "A"
$ws.Cells.Item(131, 12).Value = 3.75
$ws.Cells.Item(131, 13).Value = 3.5
$ws.Cells.Item(131, 14).Value = 1.8
$ws.Cells.Item(131, 15).Value = 4.5
$ws.Cells.Item(131, 16).Value = 3.6
$ws.Cells.Item(131, 17).Value = 1.65
$ws.Cells.Item(131, 18).Value = 0.75
$ws.Cells.Item(131, 19).Value = 1.925
$ws.Cells.Item(131, 20).Value = 1.875
$ws.Cells.Item(131, 21).Value = 2.25
$ws.Cells.Item(131, 22).Value = 1.9
$ws.Cells.Item(131, 23).Value = 1.9
$ws.Cells.Item(131, 24).Value = -1
$ws.Cells.Item(131, 25).Value = -1
$ws.Cells.Item(131, 26).Value = 0.6499999999999999
$ws.Cells.Item(131, 27).Value = -1
$ws.Cells.Item(131, 28).Value = 0.875
$ws.Cells.Item(131, 29).Value = 0.8999999999999999
$ws.Cells.Item(131, 30).Value = -1

# Row 132
$ws.Cells.Item(132, 2).Value = 6834846
$ws.Cells.Item(132, 5).Value = "NK Dugopolje"
$ws.Cells.Item(132, 6).Value = "Orijent"
$ws.Cells.Item(132, 7).Value = 1
$ws.Cells.Item(132, 8).Value = 2
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 1
$ws.Cells.Item(132, 11).Value = "A"
$ws.Cells.Item(132, 12).Value = 1.75
$ws.Cells.Item(132, 13).Value = 3.5
$ws.Cells.Item(132, 14).Value = 4
$ws.Cells.Item(132, 15).Value = 2
$ws.Cells.Item(132, 16).Value = 3
$ws.Cells.Item(132, 17).Value = 3.6
$ws.Cells.Item(132, 18).Value = -0.25
$ws.Cells.Item(132, 19).Value = 1.725
$ws.Cells.Item(132, 20).Value = 1.975
$ws.Cells.Item(132, 21).Value = 2.25
$ws.Cells.Item(132, 22).Value = 1.9
$ws.Cells.Item(132, 23).Value = 1.9
$ws.Cells.Item(132, 24).Value = -1
$ws.Cells.Item(132, 25).Value = -1
$ws.Cells.Item(132, 26).Value = 2.6
$ws.Cells.Item(132, 27).Value = -1
$ws.Cells.Item(132, 28).Value = 0.9750000000000001
$ws.Cells.Item(132, 29).Value = 0.8999999999999999
$ws.Cells.Item(132, 30).Value = -1

# Row 142
$ws.Cells.Item(142, 2).Value = 7977238
$ws.Cells.Item(142, 5).Value = "NK Jarun"
$ws.Cells.Item(142, 6).Value = "Bijelo Brdo"
$ws.Cells.Item(142, 7).Value = 1
$ws.Cells.Item(142, 8).Value = 3
$ws.Cells.Item(142, 9).Value = 1
$ws.Cells.Item(142, 10).Value = 2
$ws.Cells.Item(142, 11).Value = "A"
$ws.Cells.Item(142, 12).Value = 1.909
$ws.Cells.Item(142, 13).Value = 3.3
$ws.Cells.Item(142, 14).Value = 3.5
$ws.Cells.Item(142, 15).Value = 1.909
$ws.Cells.Item(142, 16).Value = 3.4
$ws.Cells.Item(142, 17).Value = 3.5
$ws.Cells.Item(142, 18).Value = -0.5
$ws.Cells.Item(142, 19).Value = 1.95
$ws.Cells.Item(142, 20).Value = 1.85
$ws.Cells.Item(142, 21).Value = 2.25
$ws.Cells.Item(142, 22).Value = 1.75
$ws.Cells.Item(142, 23).Value = 2.05
$ws.Cells.Item(142, 24).Value = -1
$ws.Cells.Item(142, 25).Value = -1
$ws.Cells.Item(142, 26).Value = 2.5
$ws.Cells.Item(142, 27).Value = -1
$ws.Cells.Item(142, 28).Value = 0.8500000000000001
$ws.Cells.Item(142, 29).Value = 0.75
$ws.Cells.Item(142, 30).Value = -1

# Row 143
$ws.Cells.Item(143, 2).Value = 7977239
$ws.Cells.Item(143, 5).Value = "NK Dugopolje"
$ws.Cells.Item(143, 6).Value = "NK Solin"
$ws.Cells.Item(143, 7).Value = 2
$ws.Cells.Item(143, 8).Value = 2
$ws.Cells.Item(143, 9).Value = 0
$ws.Cells.Item(143, 10).Value = 1
$ws.Cells.Item(143, 11).Value = "D"
$ws.Cells.Item(143, 12).Value = 2
$ws.Cells.Item(143, 13).Value = 3.2
$ws.Cells.Item(143, 14).Value = 3.3
$ws.Cells.Item(143, 15).Value = 2.3
$ws.Cells.Item(143, 16).Value = 3.2
$ws.Cells.Item(143, 17).Value = 2.7
$ws.Cells.Item(143, 18).Value = -0.25
$ws.Cells.Item(143, 19).Value = 1.975
$ws.Cells.Item(143, 20).Value = 1.725
$ws.Cells.Item(143, 21).Value = 2.25
$ws.Cells.Item(143, 22).Value = 1.9
$ws.Cells.Item(143, 23).Value = 1.9
$ws.Cells.Item(143, 24).Value = -1
$ws.Cells.Item(143, 25).Value = 2.2
$ws.Cells.Item(143, 26).Value = -1
$ws.Cells.Item(143, 27).Value = -0.5
$ws.Cells.Item(143, 28).Value = 0.3625
$ws.Cells.Item(143, 29).Value = 0.8999999999999999
$ws.Cells.Item(143, 30).Value = -1

# Row 147
$ws.Cells.Item(147, 2).Value = 7977243
$ws.Cells.Item(147, 5).Value = "Bijelo Brdo"
$ws.Cells.Item(147, 6).Value = "NK Sesvete"
$ws.Cells.Item(147, 7).Value = 2
$ws.Cells.Item(147, 8).Value = 1
$ws.Cells.Item(147, 9).Value = 0
$ws.Cells.Item(147, 10).Value = 1
$ws.Cells.Item(147, 11).Value = "H"
$ws.Cells.Item(147, 12).Value = 2.4
$ws.Cells.Item(147, 13).Value = 3.4
$ws.Cells.Item(147, 14).Value = 2.5
$ws.Cells.Item(147, 15).Value = 1.909
$ws.Cells.Item(147, 16).Value = 3.2
$ws.Cells.Item(147, 17).Value = 4
$ws.Cells.Item(147, 18).Value = -0.5
$ws.Cells.Item(147, 19).Value = 1.9
$ws.Cells.Item(147, 20).Value = 1.9
$ws.Cells.Item(147, 21).Value = 2.25
$ws.Cells.Item(147, 22).Value = 1.975
$ws.Cells.Item(147, 23).Value = 1.725
$ws.Cells.Item(147, 24).Value = 0.909
$ws.Cells.Item(147, 25).Value = -1
$ws.Cells.Item(147, 26).Value = -1
$ws.Cells.Item(147, 27).Value = 0.8999999999999999
$ws.Cells.Item(147, 28).Value = -1
$ws.Cells.Item(147, 29).Value = 0.9750000000000001
$ws.Cells.Item(147, 30).Value = -1

# Row 148
$ws.Cells.Item(148, 2).Value = 7977245
$ws.Cells.Item(148, 5).Value = "NK Dugopolje"
$ws.Cells.Item(148, 6).Value = "NK Dubrava Zagreb"
$ws.Cells.Item(148, 7).Value = 0
$ws.Cells.Item(148, 8).Value = 0
$ws.Cells.Item(148, 9).Value = 0
$ws.Cells.Item(148, 10).Value = 0
$ws.Cells.Item(148, 11).Value = "D"
$ws.Cells.Item(148, 12).Value = 2
$ws.Cells.Item(148, 13).Value = 3.1
$ws.Cells.Item(148, 14).Value = 3.4
$ws.Cells.Item(148, 15).Value = 2.7
$ws.Cells.Item(148, 16).Value = 3
$ws.Cells.Item(148, 17).Value = 2.55
$ws.Cells.Item(148, 18).Value = 0
$ws.Cells.Item(148, 19).Value = 1.95
$ws.Cells.Item(148, 20).Value = 1.85
$ws.Cells.Item(148, 21).Value = 2.25
$ws.Cells.Item(148, 22).Value = 1.925
$ws.Cells.Item(148, 23).Value = 1.875
$ws.Cells.Item(148, 24).Value = -1
$ws.Cells.Item(148, 25).Value = 2
$ws.Cells.Item(148, 26).Value = -1
$ws.Cells.Item(148, 27).Value = 0
$ws.Cells.Item(148, 28).Value = 0
$ws.Cells.Item(148, 29).Value = -1
$ws.Cells.Item(148, 30).Value = 0.875
